# Edit script for ZBP_12_obavy_ztrata_prace.xlsx style update
# Adds a new date column (22. 2. 2022) to both sheets, fixes the AM1/AL1 header
# typo (2021 -> 2022), updates several re-based percentage/count values, and
# refreshes the footer "aktualizace" text.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")
$ws2 = $wb.Worksheets.Item("pocetR")

# --- Sheet "data": new column AN, cloned formatting from header cell AM1 ---
$ws1.Range("AM1").Copy($ws1.Range("AN1"))

$ws1.Range("AM1").Value = "25. 1. 2022"
$ws1.Range("AN1").Value = "22. 2. 2022"
$ws1.Range("AN2").Value = 0.5
$ws1.Range("AN3").Value = 0.32
$ws1.Range("AN4").Value = 0.18
$ws1.Range("AN5").Value = 0.27
$ws1.Range("AM6").Value = 0.27
$ws1.Range("AN6").Value = 0.32
$ws1.Range("AM7").Value = 0.35
$ws1.Range("AN7").Value = 0.41
$ws1.Range("AM8").Value = 0.54
$ws1.Range("AN8").Value = 0.53
$ws1.Range("AM9").Value = 0.33
$ws1.Range("AN9").Value = 0.31
$ws1.Range("AN10").Value = 0.16
$ws1.Range("AM11").Value = 0.54
$ws1.Range("AN11").Value = 0.51
$ws1.Range("AM12").Value = 0.33
$ws1.Range("AN12").Value = 0.33
$ws1.Range("AN13").Value = 0.16
$ws1.Range("AM14").Value = 0.51
$ws1.Range("AN14").Value = 0.5
$ws1.Range("AM15").Value = 0.23
$ws1.Range("AN15").Value = 0.3
$ws1.Range("AN16").Value = 0.2
$ws1.Range("AM17").Value = 0.55
$ws1.Range("AN17").Value = 0.51
$ws1.Range("AM18").Value = 0.32
$ws1.Range("AN18").Value = 0.32
$ws1.Range("AN19").Value = 0.17
$ws1.Range("AN20").Value = 0.5
$ws1.Range("AN21").Value = 0.31
$ws1.Range("AN22").Value = 0.19
$ws1.Range("AM23").Value = 0.4
$ws1.Range("AN23").Value = 0.43
$ws1.Range("AM24").Value = 0.44
$ws1.Range("AN24").Value = 0.31
$ws1.Range("AM25").Value = 0.16
$ws1.Range("AN25").Value = 0.26
$ws1.Range("AN26").Value = 0.35
$ws1.Range("AN27").Value = 0.38
$ws1.Range("AN28").Value = 0.27
$ws1.Range("AM29").Value = 0.55
$ws1.Range("AN29").Value = 0.59
$ws1.Range("AM30").Value = 0.32
$ws1.Range("AN30").Value = 0.29
$ws1.Range("AM31").Value = 0.13
$ws1.Range("AN31").Value = 0.12
$ws1.Range("AN32").Value = 0.6
$ws1.Range("AN33").Value = 0.26
$ws1.Range("AN34").Value = 0.14
$ws1.Range("AM35").Value = 0.43
$ws1.Range("AN35").Value = 0.39
$ws1.Range("AM36").Value = 0.37
$ws1.Range("AN36").Value = 0.38
$ws1.Range("AM37").Value = 0.2
$ws1.Range("AN37").Value = 0.23
$ws1.Range("AN38").Value = 0.42
$ws1.Range("AM39").Value = 0.37
$ws1.Range("AN39").Value = 0.35
$ws1.Range("AM40").Value = 0.18
$ws1.Range("AN40").Value = 0.23
$ws1.Range("AM41").Value = 0.59
$ws1.Range("AN41").Value = 0.55
$ws1.Range("AM42").Value = 0.29
$ws1.Range("AN42").Value = 0.29
$ws1.Range("AN43").Value = 0.16
$ws1.Range("AN44").Value = 0.62
$ws1.Range("AN45").Value = 0.29
$ws1.Range("AN46").Value = 0.09
$ws1.Range("AM47").Value = 0.52
$ws1.Range("AN47").Value = 0.49
$ws1.Range("AN48").Value = 0.31
$ws1.Range("AM49").Value = 0.17
$ws1.Range("AN49").Value = 0.2
$ws1.Range("AM50").Value = 0.66
$ws1.Range("AN50").Value = 0.65
$ws1.Range("AM51").Value = 0.3
$ws1.Range("AN51").Value = 0.29
$ws1.Range("AM52").Value = 0.04
$ws1.Range("AN52").Value = 0.06
$ws1.Range("AN53").Value = 0.45
$ws1.Range("AN54").Value = 0.33
$ws1.Range("AN55").Value = 0.22
$ws1.Range("AM56").Value = 0.59
$ws1.Range("AN56").Value = 0.54
$ws1.Range("AM57").Value = 0.27
$ws1.Range("AN57").Value = 0.28
$ws1.Range("AM58").Value = 0.14
$ws1.Range("AN58").Value = 0.18
$ws1.Range("AN59").Value = 0.62
$ws1.Range("AM60").Value = 0.23
$ws1.Range("AN60").Value = 0.29
$ws1.Range("AM61").Value = 0.11
$ws1.Range("AN61").Value = 0.09
$ws1.Range("A62").Value = "Život během pandemie, Obavy ze ztráty práce, % respondentů celkově a ve skupinách, aktualizace 2. 3. 2022"

# --- Sheet "pocetR": new column AM, cloned formatting from header cell AL1 ---
$ws2.Range("AL1").Copy($ws2.Range("AM1"))

$ws2.Range("AL1").Value = "25. 1. 2022"
$ws2.Range("AM1").Value = "22. 2. 2022"
$ws2.Range("AL2").Value = 1029
$ws2.Range("AM2").Value = 1001
$ws2.Range("AL3").Value = 99
$ws2.Range("AM3").Value = 90
$ws2.Range("AL4").Value = 930
$ws2.Range("AM4").Value = 911
$ws2.Range("AL5").Value = 790
$ws2.Range("AM5").Value = 768
$ws2.Range("AL6").Value = 156
$ws2.Range("AM6").Value = 150
$ws2.Range("AM7").Value = 9
$ws2.Range("AL8").Value = 77
$ws2.Range("AM8").Value = 74
$ws2.Range("AL9").Value = 757
$ws2.Range("AM9").Value = 734
$ws2.Range("AL10").Value = 144
$ws2.Range("AM10").Value = 136
$ws2.Range("AL11").Value = 65
$ws2.Range("AM11").Value = 65
$ws2.Range("AL12").Value = 63
$ws2.Range("AM12").Value = 66
$ws2.Range("AL13").Value = 366
$ws2.Range("AM13").Value = 361
$ws2.Range("AL14").Value = 413
$ws2.Range("AM14").Value = 400
$ws2.Range("AL15").Value = 250
$ws2.Range("AM15").Value = 240
$ws2.Range("AL16").Value = 116
$ws2.Range("AM16").Value = 110
$ws2.Range("AL17").Value = 315
$ws2.Range("AM17").Value = 300
$ws2.Range("AL18").Value = 349
$ws2.Range("AM18").Value = 335
$ws2.Range("AL19").Value = 154
$ws2.Range("AM19").Value = 157
$ws2.Range("AL20").Value = 284
$ws2.Range("AM20").Value = 268
$ws2.Range("AL21").Value = 93
$ws2.Range("AM21").Value = 92
$ws2.Range("AL22").Value = 265
$ws2.Range("AM22").Value = 273
$ws2.Range("AL23").Value = 146
$ws2.Range("AM23").Value = 133
$ws2.Range("AL24").Value = 88
$ws2.Range("AM24").Value = 87
$ws2.Range("A25").Value = "Život během pandemie, Obavy ze ztráty práce, velikost dotázaného souboru celkově a ve skupinách, aktualizace 2. 3. 2022"
$ws2.Range("AL25").Copy($ws2.Range("AM25"))
